$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.727236390113831
$ws.Range("B1").Value = 2.640979766845703
$ws.Range("C1").Value = 2.813587427139282
$ws.Range("D1").Value = 3.17207932472229
$ws.Range("E1").Value = 2.956655502319336
